# Refresh the crypto price/volume table with the latest scrape results.
# Two rows (Chainlink/Polkadot, Maker/Stacks, InjectiveProtocol/dogwifhat) swapped
# rank position in this run, so Coin/Link are rewritten for those rows too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='66.314.41'; E='  -1.38%  ' }
    @{ Row=3; D='3.314.77'; E='  +1.14%  ' }
    @{ Row=4; E='  -0.21%  ' }
    @{ Row=5; D='575.23'; E='  -0.35%  ' }
    @{ Row=6; D='180.88'; E='  -2.41%  ' }
    @{ Row=7; D='0.629'; E='  +5.84%  ' }
    @{ Row=8; E='  -0.12%  ' }
    @{ Row=9; E='  -0.94%  ' }
    @{ Row=10; D='6.70'; E='  +0.32%  ' }
    @{ Row=11; D='0.404'; E='  -0.87%  ' }
    @{ Row=12; D='3.898.44'; E='  +0.95%  ' }
    @{ Row=13; E='  -3.24%  ' }
    @{ Row=14; D='26.70'; E='  -2.16%  ' }
    @{ Row=15; D='66.362.22'; E='  -1.78%  ' }
    @{ Row=16; E='  -0.65%  ' }
    @{ Row=17; D='3.296.52'; E='  -0.72%  ' }
    @{ Row=18; D='442.03'; E='  +1.01%  ' }
    @{ Row=19; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.65'; E='  -0.82%  ' }
    @{ Row=20; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='13.50'; E='  -0.45%  ' }
    @{ Row=21; D='7.54'; E='  -2.25%  ' }
    @{ Row=22; D='73.28'; E='  -1.17%  ' }
    @{ Row=23; D='0.998'; E='  +0.05%  ' }
    @{ Row=24; D='0.519'; E='  +1.24%  ' }
    @{ Row=25; D='3.468.71'; E='  +0.57%  ' }
    @{ Row=26; E='  -1.61%  ' }
    @{ Row=27; D='0.195'; E='  +3.79%  ' }
    @{ Row=28; D='9.03'; E='  -0.13%  ' }
    @{ Row=29; D='0.999'; E='  -0.42%  ' }
    @{ Row=30; E='  -0.93%  ' }
    @{ Row=31; D='22.70'; E='  -0.20%  ' }
    @{ Row=32; E='  +0.05%  ' }
    @{ Row=33; D='6.75'; E='  -0.20%  ' }
    @{ Row=34; E='  -2.38%  ' }
    @{ Row=35; E='  -1.88%  ' }
    @{ Row=36; E='  -3.01%  ' }
    @{ Row=37; D='159.49'; E='  -2.01%  ' }
    @{ Row=38; D='27.35'; E='  +1.95%  ' }
    @{ Row=39; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.82'; E='  -1.13%  ' }
    @{ Row=40; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='2.835.02'; E='  +4.64%  ' }
    @{ Row=41; D='0.788'; E='  +0.57%  ' }
    @{ Row=42; E='  -1.24%  ' }
    @{ Row=43; D='40.63'; E='  +1.01%  ' }
    @{ Row=44; D='6.17'; E='  -2.86%  ' }
    @{ Row=45; D='0.0667'; E='  -0.66%  ' }
    @{ Row=46; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='23.97'; E='  -2.21%  ' }
    @{ Row=47; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='2.33'; E='  -1.98%  ' }
    @{ Row=48; D='325.16'; E='  +0.00%  ' }
    @{ Row=49; E='  -0.07%  ' }
    @{ Row=50; D='0.103'; E='  +3.59%  ' }
    @{ Row=51; E='  -0.83%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey('D')) {
        # Prices are stored as text (some contain thousands separators like
        # '66.314.41' that aren't valid numbers); force text so Excel doesn't
        # reinterpret plain decimals (e.g. '1.00' -> 1) and drop formatting.
        $ws.Range("D$r").NumberFormat = '@'
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.ContainsKey('E')) { $ws.Range("E$r").Value = $u.E }
}
